$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of (row, col) -> new cell text, 1-indexed as in the Word object model.
# Only the five data rows (1, 5, 9, 13, 17) of this 20x5 table hold text;
# the other rows are blank spacer rows and are left untouched.
$updates = @{
    "1,1"  = "25÷9=2, 7"
    "1,2"  = "98÷9=10, 8"
    "1,3"  = "56÷5=11, 1"
    "1,4"  = "54÷8=6, 6"
    "1,5"  = "57÷4=14, 1"

    "5,1"  = "77÷4=19, 1"
    "5,2"  = "69÷4=17, 1"
    "5,3"  = "10÷5=2, 0"
    "5,4"  = "69÷4=17, 1"
    "5,5"  = "64÷7=9, 1"

    "9,1"  = "27÷2=13, 1"
    "9,2"  = "58÷4=14, 2"
    "9,3"  = "44÷4=11, 0"
    "9,4"  = "23÷5=4, 3"
    "9,5"  = "52÷9=5, 7"

    "13,1" = "27÷7=3, 6"
    "13,2" = "14÷4=3, 2"
    "13,3" = "83÷2=41, 1"
    "13,4" = "46÷4=11, 2"
    "13,5" = "89÷8=11, 1"

    "17,1" = "67÷4=16, 3"
    "17,2" = "83÷9=9, 2"
    "17,3" = "69÷2=34, 1"
    "17,4" = "32÷5=6, 2"
    "17,5" = "10÷3=3, 1"
}

foreach ($key in $updates.Keys) {
    $parts = $key.Split(",")
    $row = [int]$parts[0]
    $col = [int]$parts[1]
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    # Trim the trailing end-of-cell marker so only the visible text is replaced.
    $rng.End = $rng.End - 1
    $rng.Text = $updates[$key]
}
